$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 SVC(C=1, class_weight='balanced', kernel='sigmoid',
                     random_state=42))])"
$ws.Range("B2").Value = 0.7549999999999999
$ws.Range("C2").Value = "{'selector': None, 'scaler': None, 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 1}"
$ws.Range("D2").Value = 0.7485513576894436
$ws.Range("E2").Value = 0.6367576740389239
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.6921155419927071
$ws.Range("H2").Value = 0.583012566137566
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.8741710992907802
$ws.Range("K2").Value = 0.7573611111111112
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]"
$ws.Range("N2").Value = "[0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0]"
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model', SVC(C=5, kernel='sigmoid', random_state=42))])"
$ws.Range("B3").Value = 0.76
$ws.Range("C3").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': MinMaxScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 5}"
$ws.Range("D3").Value = 0.7473808081477757
$ws.Range("E3").Value = 0.6282883516946016
$ws.Range("F3").Value = 0.7999999999999999
$ws.Range("G3").Value = 0.6858006110630782
$ws.Range("H3").Value = 0.6013684689153439
$ws.Range("I3").Value = 0.7368421052631579
$ws.Range("J3").Value = 0.8748670212765959
$ws.Range("K3").Value = 0.728576388888889
$ws.Range("L3").Value = 0.875
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]"
$ws.Range("N3").Value = "[1 1 1 1 1 1 1 0 1 1 1 1 0 1 1 1 0 1 1 1 0 1 1 0]"
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 SVC(C=3, class_weight='balanced', kernel='poly',
                     random_state=42))])"
$ws.Range("B4").Value = 0.7041647241647242
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__kernel': 'poly', 'model__class_weight': 'balanced', 'model__C': 3}"
$ws.Range("D4").Value = 0.7426144602553902
$ws.Range("E4").Value = 0.5889400356587856
$ws.Range("F4").Value = 0.7567567567567567
$ws.Range("G4").Value = 0.6821270793350787
$ws.Range("H4").Value = 0.5269593253968254
$ws.Range("I4").Value = 0.7777777777777778
$ws.Range("J4").Value = 0.8713935185185185
$ws.Range("K4").Value = 0.7461249999999999
$ws.Range("L4").Value = 0.7368421052631579
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]"
$ws.Range("N4").Value = "[1 1 1 1 1 1 0 1 1 1 1 0 1 1 1 1 0 0 1 0 1 0 1 1]"
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),
                ('model', SVC(C=1, kernel='sigmoid', random_state=42))])"
$ws.Range("B5").Value = 0.7761172161172161
$ws.Range("C5").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 1}"
$ws.Range("D5").Value = 0.6106929038424203
$ws.Range("E5").Value = 0.5198868885743885
$ws.Range("F5").Value = 0.7777777777777778
$ws.Range("G5").Value = 0.5796088135739609
$ws.Range("H5").Value = 0.4883112599206349
$ws.Range("I5").Value = 0.6363636363636364
$ws.Range("J5").Value = 0.689655612244898
$ws.Range("K5").Value = 0.6030902777777778
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]"
$ws.Range("N5").Value = "[1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 1]"
$ws.Range("O5").Value = 99

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 SVC(C=1, class_weight='balanced', kernel='sigmoid',
                     random_state=42))])"
$ws.Range("B6").Value = 0.76
$ws.Range("C6").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': MinMaxScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 1}"
$ws.Range("D6").Value = 0.7986116345103841
$ws.Range("E6").Value = 0.6703988546176046
$ws.Range("F6").Value = 0.6285714285714286
$ws.Range("G6").Value = 0.7335917955453835
$ws.Range("H6").Value = 0.6031261574074073
$ws.Range("I6").Value = 0.4583333333333333
$ws.Range("J6").Value = 0.9093309294871794
$ws.Range("K6").Value = 0.7994791666666666
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]"
$ws.Range("N6").Value = "[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("O6").Value = 89
